$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "zipcode" header in column O (18th shared string O1)
$ws.Range("O1").Value = "zipcode"

# Fill the zipcode value (98058) for every data row (rows 2 through 456)
$ws.Range("O2:O456").Value = 98058

# Move active cell selection to O1 as in the final workbook state
$ws.Range("O1").Select() | Out-Null
